# Bug 3111 final changes for release
#
# 1. Make a backup copy of the existing "V1.2.0 1_26_2015" sheet (Excel
#    auto-names the copy "V1.2.0 1_26_2015 (2)" and places it before the
#    original tab).
# 2. Rename the original sheet to "V1.2.1 6_9_2015" for the new release.
# 3. Update the pass/fail ("x") marks on the renamed sheet to reflect the
#    new release's test results: simulator runs (32/64-bit sim) passed,
#    device runs did not (yet), so clear the old all-green marks and set
#    only the simulator columns.
# 4. Leave the new "V1.2.1 6_9_2015" sheet active/selected, matching the
#    last place the author was working (F21).

$wb = $excel.ActiveWorkbook

$orig = $wb.Worksheets.Item("V1.2.0 1_26_2015")

# Step 1: copy the sheet to preserve the old release's data as a backup.
# (Copying in place, before itself, matches Excel's "Move or Copy..." with
# "Create a copy" checked and the copy dropped ahead of the source tab.)
$orig.Copy($orig)

# Step 2: rename the original (now second) tab for the new release. Re-fetch
# it by its still-original name rather than reusing the $orig handle, since
# that handle now follows the freshly created copy.
$toRename = $wb.Worksheets.Item("V1.2.0 1_26_2015")
$toRename.Name = "V1.2.1 6_9_2015"

# Step 3: update the test-status cells on the new release sheet.
$ws = $wb.Worksheets.Item("V1.2.1 6_9_2015")

$clearCells = @(
    "D5","E5","F5","G5",
    "C6","E6","G6",
    "C7","E7",
    "C8","E8","G8",
    "C9","E9","G9",
    "C10","E10","G10",
    "C11","E11","G11",
    "C12","E12","G12",
    "C13","E13","G13",
    "C14","E14","G14",
    "C15","E15","G15",
    "C16","E16","G16",
    "C17","E17","G17",
    "C19",
    "C20",
    "C21"
)
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

$setCells = @("C18","D19","F19","D20","F20")
foreach ($addr in $setCells) {
    $ws.Range($addr).Value = "x"
}

# Step 4: leave the new release sheet active, selection parked at F21.
$ws.Activate()
$ws.Range("F21").Select()
